# Updated symbol list on Mon Jan 30 15:36:17 UTC 2023 with GitHub Actions
#
# Applies the refreshed Price (column D) and Volume(1h) (column E) values
# from the coinranking.com scrape. Values are written as literal text
# (matching the workbook's existing inlineStr storage) rather than as
# numbers/percentages, so Excel's automatic numeric/percent type-detection
# and number-format changes are explicitly undone afterwards.

function Set-CellText($ws, $addr, $text) {
    # Prefix with an apostrophe so Excel treats the input as literal text
    # instead of auto-converting it to a number/percentage.
    $ws.Range($addr).Value = "'" + $text
    # Excel applies a "quote prefix" text style when coercing to text;
    # restore the original (default/"Normal") cell style so formatting
    # stays exactly as it was before the edit.
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "308.71"
Set-CellText $ws "E2" "-3.23%"
Set-CellText $ws "D3" "37.90"
Set-CellText $ws "E3" "-4.60%"
Set-CellText $ws "D4" "5.076"
Set-CellText $ws "E4" "-1.07%"
Set-CellText $ws "D5" "0.07880"
Set-CellText $ws "E5" "-3.84%"
Set-CellText $ws "D6" "1.965"
Set-CellText $ws "E6" "-2.69%"
Set-CellText $ws "D7" "4.353"
Set-CellText $ws "E7" "1.79%"
Set-CellText $ws "E8" "0.11%"
Set-CellText $ws "D9" "3.035"
Set-CellText $ws "E9" "-4.15%"
Set-CellText $ws "D10" "0.9304"
Set-CellText $ws "E10" "-0.55%"
Set-CellText $ws "D11" "0.1306"
Set-CellText $ws "E11" "-7.60%"
Set-CellText $ws "D12" "0.1998"
Set-CellText $ws "E12" "0.11%"
Set-CellText $ws "D13" "0.08784"
Set-CellText $ws "E13" "-3.27%"
Set-CellText $ws "D14" "0.03429"
Set-CellText $ws "E14" "-2.65%"
Set-CellText $ws "D15" "0.09724"
Set-CellText $ws "E15" "-0.76%"
Set-CellText $ws "D16" "0.001389"
Set-CellText $ws "E16" "-0.68%"
Set-CellText $ws "D17" "0.005900"
Set-CellText $ws "E17" "-5.53%"
Set-CellText $ws "E18" "1,777.17%"
Set-CellText $ws "D19" "3.590"
Set-CellText $ws "E19" "-2.00%"
Set-CellText $ws "D20" "0.3474"
Set-CellText $ws "E20" "0.36%"
Set-CellText $ws "E21" "1.56%"
Set-CellText $ws "D22" "5.006"
Set-CellText $ws "E22" "2.16%"
Set-CellText $ws "D23" "0.2489"
Set-CellText $ws "E23" "1.45%"
Set-CellText $ws "D24" "0.04323"
Set-CellText $ws "E24" "-0.13%"
Set-CellText $ws "D25" "0.001219"
Set-CellText $ws "E25" "-0.46%"
Set-CellText $ws "D26" "0.004612"
Set-CellText $ws "E26" "-3.51%"
Set-CellText $ws "D27" "0.0001352"
Set-CellText $ws "E27" "4.00%"
Set-CellText $ws "D39" "0.02285"
Set-CellText $ws "E39" "2.54%"
Set-CellText $ws "D40" "0.05036"
Set-CellText $ws "E40" "-4.80%"
Set-CellText $ws "D41" "0.007513"
Set-CellText $ws "D42" "0.009861"
Set-CellText $ws "E42" "-0.76%"
Set-CellText $ws "D43" "0.1358"
Set-CellText $ws "E43" "-1.40%"
Set-CellText $ws "D44" "0.002043"
Set-CellText $ws "E44" "-4.97%"
Set-CellText $ws "D45" "0.008767"
Set-CellText $ws "E45" "-11.20%"
Set-CellText $ws "D46" "0.00006585"
Set-CellText $ws "E46" "1.11%"
Set-CellText $ws "E47" "0.05%"
Set-CellText $ws "D48" "0.002998"
Set-CellText $ws "E48" "8.26%"
Set-CellText $ws "D50" "0.00002100"
Set-CellText $ws "E50" "0.05%"
Set-CellText $ws "D51" "0.0002000"
Set-CellText $ws "E51" "0.05%"
